$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 195223.89
$ws.Range("J17").Value = 195223.89
$ws.Range("L17").Value = 585671.67
$ws.Range("N17").Value = -586007.67
$ws.Range("H113").Value = 12999.5
$ws.Range("I113").Value = 21999
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 21999
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -18745
$ws.Range("N113").Value = -10508
$ws.Range("H136").Value = 40363.637
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H137").Value = 3850.7778
$ws.Range("I137").Value = 4652.4546
$ws.Range("J137").Value = 2591
$ws.Range("K137").Value = 13957.3638
$ws.Range("L137").Value = 7773
$ws.Range("M137").Value = -11407.3638
$ws.Range("N137").Value = -12873

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1792.3684
$ws.Range("I2").Value = 1634.3572
$ws.Range("J2").Value = 2234.8
$ws.Range("K2").Value = 1634.3572
$ws.Range("L2").Value = 2234.8
$ws.Range("M2").Value = -1521.3572
$ws.Range("N2").Value = -2460.8
$ws.Range("H61").Value = 7152
$ws.Range("I61").Value = 6739.913
$ws.Range("J61").Value = 7430.7646
$ws.Range("K61").Value = 6739.913
$ws.Range("L61").Value = 7430.7646
$ws.Range("M61").Value = -6527.913
$ws.Range("N61").Value = -7854.7646
$ws.Range("H74").Value = 2816.25
$ws.Range("I74").Value = 848.61536
$ws.Range("K74").Value = 848.61536
$ws.Range("M74").Value = 25.38463999999999
$ws.Range("H77").Value = 2816.25
$ws.Range("I77").Value = 848.61536
$ws.Range("K77").Value = 4243.0768
$ws.Range("M77").Value = 124.9232000000002
$ws.Range("H116").Value = 1792.3684
$ws.Range("I116").Value = 1634.3572
$ws.Range("J116").Value = 2234.8
$ws.Range("K116").Value = 1634.3572
$ws.Range("L116").Value = 2234.8
$ws.Range("M116").Value = 659.6428000000001
$ws.Range("N116").Value = -6822.8
$ws.Range("H136").Value = 7152
$ws.Range("I136").Value = 6739.913
$ws.Range("J136").Value = 7430.7646
$ws.Range("K136").Value = 20219.739
$ws.Range("L136").Value = 22292.2938
$ws.Range("M136").Value = -17669.739
$ws.Range("N136").Value = -27392.2938
$ws.Range("H141").Value = 94429
$ws.Range("J141").Value = 94429
$ws.Range("L141").Value = 94429
$ws.Range("N141").Value = -104789

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1792.3684
$ws.Range("I3").Value = 1634.3572
$ws.Range("J3").Value = 2234.8
$ws.Range("K3").Value = 1634.3572
$ws.Range("L3").Value = 2234.8
$ws.Range("M3").Value = -1520.3572
$ws.Range("N3").Value = -2462.8
$ws.Range("H20").Value = 1598.1666
$ws.Range("I20").Value = 1395.3334
$ws.Range("J20").Value = 1902.4166
$ws.Range("K20").Value = 1395.3334
$ws.Range("L20").Value = 1902.4166
$ws.Range("M20").Value = -1148.3334
$ws.Range("N20").Value = -2396.4166
$ws.Range("H94").Value = 2718.6956
$ws.Range("J94").Value = 13264.75
$ws.Range("L94").Value = 13264.75
$ws.Range("N94").Value = -14166.75
$ws.Range("H107").Value = 3092.1482
$ws.Range("I107").Value = 2578.4707
$ws.Range("K107").Value = 2578.4707
$ws.Range("M107").Value = -658.4706999999999
$ws.Range("H132").Value = 99999.664
$ws.Range("J132").Value = 99999.664
$ws.Range("L132").Value = 99999.664
$ws.Range("N132").Value = -110119.664
$ws.Range("H134").Value = 6399.1514
$ws.Range("I134").Value = 3254.9333
$ws.Range("K134").Value = 9764.7999
$ws.Range("M134").Value = -7229.7999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4128.731
$ws.Range("I31").Value = 2552.4
$ws.Range("K31").Value = 2552.4
$ws.Range("M31").Value = -2257.4
$ws.Range("H34").Value = 4128.731
$ws.Range("I34").Value = 2552.4
$ws.Range("K34").Value = 2552.4
$ws.Range("M34").Value = -2350.4
$ws.Range("H122").Value = 4449.3335
$ws.Range("I122").Value = 4549.6665
$ws.Range("J122").Value = 4248.6665
$ws.Range("K122").Value = 13648.9995
$ws.Range("L122").Value = 12745.9995
$ws.Range("M122").Value = -11198.9995
$ws.Range("N122").Value = -17645.9995
$ws.Range("H123").Value = 147499.75
$ws.Range("J123").Value = 109999.664
$ws.Range("L123").Value = 109999.664
$ws.Range("N123").Value = -119799.664

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 137.03703
$ws.Range("J2").Value = 181.75
$ws.Range("L2").Value = 1090.5
$ws.Range("N2").Value = -1316.5
$ws.Range("H4").Value = 8142306.5
$ws.Range("I4").Value = 9140110
$ws.Range("K4").Value = 27420330
$ws.Range("M4").Value = -27420218
$ws.Range("H47").Value = 397.5
$ws.Range("I47").Value = 462.66666
$ws.Range("K47").Value = 1387.99998
$ws.Range("M47").Value = -956.9999800000001
$ws.Range("H121").Value = 6585.143
$ws.Range("J121").Value = 6585.143
$ws.Range("L121").Value = 19755.429
$ws.Range("N121").Value = -22375.429
$ws.Range("H129").Value = 63542.125
$ws.Range("I129").Value = 83806.914
$ws.Range("K129").Value = 251420.742
$ws.Range("M129").Value = -246420.742

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4114.5835
$ws.Range("I122").Value = 3151.7334
$ws.Range("J122").Value = 5719.3335
$ws.Range("K122").Value = 9455.200199999999
$ws.Range("L122").Value = 17158.0005
$ws.Range("M122").Value = -7005.200199999999
$ws.Range("N122").Value = -22058.0005
$ws.Range("H132").Value = 3667.0188
$ws.Range("I132").Value = 3357.6216
$ws.Range("J132").Value = 4382.5
$ws.Range("K132").Value = 10072.8648
$ws.Range("L132").Value = 13147.5
$ws.Range("M132").Value = -7542.864799999999
$ws.Range("N132").Value = -18207.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3005.9565
$ws.Range("I22").Value = 587.5
$ws.Range("J22").Value = 3236.2856
$ws.Range("K22").Value = 587.5
$ws.Range("L22").Value = 3236.2856
$ws.Range("M22").Value = -292.5
$ws.Range("N22").Value = -3826.2856
$ws.Range("H27").Value = 3005.9565
$ws.Range("I27").Value = 587.5
$ws.Range("J27").Value = 3236.2856
$ws.Range("K27").Value = 587.5
$ws.Range("L27").Value = 3236.2856
$ws.Range("M27").Value = -480.5
$ws.Range("N27").Value = -3450.2856
$ws.Range("H115").Value = 89997.5
$ws.Range("J115").Value = 89997.5
$ws.Range("L115").Value = 89997.5
$ws.Range("N115").Value = -92347.5
$ws.Range("H125").Value = 105000
$ws.Range("J125").Value = 105000
$ws.Range("L125").Value = 105000
$ws.Range("N125").Value = -114840
$ws.Range("H132").Value = 3094.1633
$ws.Range("I132").Value = 3099.2666
$ws.Range("J132").Value = 3086.1052
$ws.Range("K132").Value = 9297.799800000001
$ws.Range("L132").Value = 9258.3156
$ws.Range("M132").Value = -6767.799800000001
$ws.Range("N132").Value = -14318.3156

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29999.092
$ws.Range("H73").Value = 29999.092
$ws.Range("H96").Value = 3285.5
$ws.Range("I96").Value = 1894
$ws.Range("K96").Value = 1894
$ws.Range("M96").Value = -521
$ws.Range("H119").Value = 98000
$ws.Range("J119").Value = 98000
$ws.Range("L119").Value = 98000
$ws.Range("N119").Value = -107676
$ws.Range("H132").Value = 4344.8887
$ws.Range("I132").Value = 4872.143
$ws.Range("K132").Value = 14616.429
$ws.Range("M132").Value = -12086.429
